{"js": "// Word JS API (Office.js) edit script.\n// Body of: async (context) => { ... }\n//\n// The \"DELUXE COFFEE\" sellers document has a 3-column schedule table\n// (FRANJA HORARIA | MESA | COMPRADOR). The edit:\n//  - swaps several buyer names between time slots,\n//  - removes the \"10:00 - 10:15\" slot,\n//  - inserts a new \"09:30 - 09:45\" slot (with buyer \"COLFRESH COFFEE\"),\n// while keeping the header row and the total row count (9 incl. header)\n// unchanged. Rewriting the whole table's `values` in one shot keeps the\n// row objects/formatting intact and sidesteps any row insert/delete\n// ordering pitfalls.\n\nconst table = context.document.body.tables.getFirstOrNullObject();\ntable.load(\"values,rowCount\");\nawait context.sync();\n\nif (table.isNullObject) {\n  throw new Error(\"Expected a table in the document body, found none.\");\n}\n\n// Build a lookup of the current FRANJA HORARIA -> COMPRADOR mapping so the\n// edit is resilient to the exact current row order.\nconst header = table.values[0];\nconst dataRows = table.values.slice(1);\n\nconst byTime = new Map();\nfor (const row of dataRows) {\n  byTime.set(row[0], row);\n}\n\n// Target buyer per time slot, after the edit (see commit message / diff):\n//   08:45 - 09:00 : ARMANDO VEL\u00c1SQUEZ -> CAF\u00c9 MOLINA\n//   09:00 - 09:15 : CAF\u00c9 MOLINA       -> BOX BRAND\n//   09:45 - 10:00 : BOX BRAND         -> FLOR A FRUTO\n//   10:15 - 10:30 : FLOR A FRUTO      -> ARMANDO VEL\u00c1SQUEZ\n//   10:30 - 10:45 : COLFRESH COFFEE   -> INTERLINK2AMERICAS\n// and the \"10:00 - 10:15\" / INTERLINK2AMERICAS row is replaced by a new\n// \"09:30 - 09:45\" / COLFRESH COFFEE row.\nconst buyerOverrides = {\n  \"08:45 - 09:00\": \"CAF\u00c9 MOLINA\",\n  \"09:00 - 09:15\": \"BOX BRAND\",\n  \"09:45 - 10:00\": \"FLOR A FRUTO\",\n  \"10:15 - 10:30\": \"ARMANDO VEL\u00c1SQUEZ\",\n  \"10:30 - 10:45\": \"INTERLINK2AMERICAS\",\n};\n\nconst newOrderOfTimes = [\n  \"08:30 - 08:45\",\n  \"08:45 - 09:00\",\n  \"09:00 - 09:15\",\n  \"09:15 - 09:30\",\n  \"09:30 - 09:45\",\n  \"09:45 - 10:00\",\n  \"10:15 - 10:30\",\n  \"10:30 - 10:45\",\n];\n\nconst newDataRows = newOrderOfTimes.map((time) => {\n  if (time === \"09:30 - 09:45\") {\n    // Brand-new row inserted by the edit.\n    return [time, \"\", \"COLFRESH COFFEE\"];\n  }\n  const existing = byTime.get(time) || [time, \"\", \"\"];\n  const mesa = existing[1] || \"\";\n  const buyer = Object.prototype.hasOwnProperty.call(buyerOverrides, time)\n    ? buyerOverrides[time]\n    : existing[2];\n  return [time, mesa, buyer];\n});\n\n// The row count before and after the edit is identical (one slot removed,\n// one added), so rewriting cell-by-cell (what the `values` setter does)\n// never needs to physically add/remove table rows.\nif (newDataRows.length + 1 !== dataRows.length + 1) {\n  throw new Error(\n    `Unexpected row-count change: before=${dataRows.length + 1}, after=${newDataRows.length + 1}`\n  );\n}\n\ntable.values = [header, ...newDataRows];\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n#\n# The \"DELUXE COFFEE\" sellers document has a 3-column schedule table\n# (FRANJA HORARIA | MESA | COMPRADOR). The edit:\n#  - swaps several buyer names between time slots,\n#  - removes the \"10:00 - 10:15\" slot,\n#  - inserts a new \"09:30 - 09:45\" slot (with buyer \"COLFRESH COFFEE\"),\n# while keeping the header row and the total row count (9 incl. header)\n# unchanged. We rewrite each data row's FRANJA HORARIA / COMPRADOR cell\n# text in place (same row count before/after) instead of physically\n# inserting/deleting table rows, which keeps per-row formatting intact\n# and avoids row-index bookkeeping.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Current time-slot -> row index (row 1 is the header row).\n$rowOfTime = @{}\nfor ($r = 2; $r -le $t.Rows.Count; $r++) {\n    $time = $t.Cell($r, 1).Range.Text.TrimEnd([char]13, [char]7)\n    $rowOfTime[$time] = $r\n}\n\n# Current COMPRADOR per time slot, read before any writes so overlapping\n# swaps (e.g. ARMANDO VEL\u00c1SQUEZ -> CAF\u00c9 MOLINA and CAF\u00c9 MOLINA -> BOX BRAND)\n# don't clobber each other.\n$buyerOfTime = @{}\nforeach ($time in $rowOfTime.Keys) {\n    $r = $rowOfTime[$time]\n    $buyerOfTime[$time] = $t.Cell($r, 3).Range.Text.TrimEnd([char]13, [char]7)\n}\n\n# Target buyer per time slot, after the edit (see commit message / diff):\n#   08:45 - 09:00 : ARMANDO VEL\u00c1SQUEZ -> CAF\u00c9 MOLINA\n#   09:00 - 09:15 : CAF\u00c9 MOLINA       -> BOX BRAND\n#   09:45 - 10:00 : BOX BRAND         -> FLOR A FRUTO\n#   10:15 - 10:30 : FLOR A FRUTO      -> ARMANDO VEL\u00c1SQUEZ\n#   10:30 - 10:45 : COLFRESH COFFEE   -> INTERLINK2AMERICAS\n# and the \"10:00 - 10:15\" / INTERLINK2AMERICAS row is replaced by a new\n# \"09:30 - 09:45\" / COLFRESH COFFEE row.\n$buyerOverrides = @{\n    \"08:45 - 09:00\" = \"CAF\u00c9 MOLINA\"\n    \"09:00 - 09:15\" = \"BOX BRAND\"\n    \"09:45 - 10:00\" = \"FLOR A FRUTO\"\n    \"10:15 - 10:30\" = \"ARMANDO VEL\u00c1SQUEZ\"\n    \"10:30 - 10:45\" = \"INTERLINK2AMERICAS\"\n}\n\n$newOrderOfTimes = @(\n    \"08:30 - 08:45\",\n    \"08:45 - 09:00\",\n    \"09:00 - 09:15\",\n    \"09:15 - 09:30\",\n    \"09:30 - 09:45\",\n    \"09:45 - 10:00\",\n    \"10:15 - 10:30\",\n    \"10:30 - 10:45\"\n)\n\n# The row count before and after the edit is identical (one slot removed,\n# one added), so rewriting cell-by-cell never needs to physically add/remove\n# table rows.\nif ($newOrderOfTimes.Length -ne $rowOfTime.Keys.Count) {\n    throw \"Unexpected row-count change: before=$($rowOfTime.Keys.Count), after=$($newOrderOfTimes.Length)\"\n}\n\nfor ($i = 0; $i -lt $newOrderOfTimes.Length; $i++) {\n    $time = $newOrderOfTimes[$i]\n    $r = $i + 2\n\n    if ($time -eq \"09:30 - 09:45\") {\n        # Brand-new row inserted by the edit.\n        $buyer = \"COLFRESH COFFEE\"\n    } elseif ($buyerOverrides.ContainsKey($time)) {\n        $buyer = $buyerOverrides[$time]\n    } else {\n        $buyer = $buyerOfTime[$time]\n    }\n\n    $t.Cell($r, 1).Range.Text = $time\n    $t.Cell($r, 3).Range.Text = $buyer\n}\n"}
